# Generate Report for Handback
#
# A new handback run completed for the "79d2a641-b11c-4d0e-b0b3-cee54dd48687"
# file, so its "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) on row 2 of both the "zh-cn"
# and "de-de" report sheets are refreshed with the new timestamps recorded
# by the report generator.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 08:37:59"
$wsZhCn.Range("H2").Value = "2016-03-20 08:38:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 08:38:03"
$wsDeDe.Range("H2").Value = "2016-03-20 08:38:25"
